$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "33.691.92"
$ws.Range("E2").Value = "  +6.88%  "

$ws.Range("D3").Value = "1.774.53"
$ws.Range("E3").Value = "  +3.83%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.556"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.53%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.277"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.83%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0661"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0921"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.14%  "

$ws.Range("D13").Value = "2.026.43"
$ws.Range("E13").Value = "  +3.70%  "

$ws.Range("D14").Value = "1.771.71"
$ws.Range("E14").Value = "  +3.73%  "

$ws.Range("E15").Value = "  +1.42%  "

$ws.Range("D16").Value = "33.659.63"
$ws.Range("E16").Value = "  +6.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "10.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.31%  "

$ws.Range("E18").Value = "  -0.92%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "249.67"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.54%  "

$ws.Range("D21").Value = "0.0₃0736"
$ws.Range("E21").Value = "  +1.49%  "

$ws.Range("E22").Value = "  +0.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.16%  "

$ws.Range("E24").Value = "  -2.44%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.15%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.95%  "

$ws.Range("E28").Value = "  +0.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.11%  "

$ws.Range("E30").Value = "  -0.09%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0514"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.28%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.54"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.83"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.88%  "

$ws.Range("D36").Value = "1.477.29"
$ws.Range("E36").Value = "  -3.27%  "

$ws.Range("E37").Value = "  +2.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.626"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.92%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0185"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.88"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.47%  "

$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.69"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.885"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.70%  "

$ws.Range("E44").Value = "  +0.87%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0510"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.94%  "

$ws.Range("D47").Value = "1.918.62"
$ws.Range("E47").Value = "  +3.71%  "

$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.69"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.92%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.87%  "
